$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post at row 472 ("好奇心旺盛なホッキョクギツネ" / curious arctic fox tweet)
# was removed from the source data. Delete that entire row; Excel shifts
# every row below it up by one and the sheet's used range shrinks from
# A1:C505 to A1:C504 automatically.
$ws.Rows(472).Delete()
